# Updated cryptos list on Wed Dec 13 21:56:53 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and fixes the OKB / Dogecoin row ordering (rows 11-12 swap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.940.97'
$ws.Range("E2").Value = '  +4.49%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.261.82'
$ws.Range("E3").Value = '  +4.07%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.32%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.13'
$ws.Range("E5").Value = '  -0.37%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.639'
$ws.Range("E6").Value = '  +2.23%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.89'
$ws.Range("E7").Value = '  +6.21%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.19%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.655'
$ws.Range("E9").Value = '  +15.30%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.11'
$ws.Range("E10").Value = '  +11.05%  '
# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0973'
$ws.Range("E11").Value = '  +4.76%  '
# Row 12
$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.68'
$ws.Range("E12").Value = '  +1.28%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.46'
$ws.Range("E13").Value = '  +6.39%  '
# Row 14
$ws.Range("E14").Value = '  +0.73%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.593.96'
$ws.Range("E15").Value = '  +3.36%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.891'
$ws.Range("E16").Value = '  +2.71%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.84'
$ws.Range("E17").Value = '  +3.69%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.252.42'
$ws.Range("E18").Value = '  +4.16%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.780.55'
$ws.Range("E19").Value = '  +4.03%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0982'
$ws.Range("E20").Value = '  +3.17%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.27'
$ws.Range("E21").Value = '  +2.43%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.61'
$ws.Range("E22").Value = '  +2.88%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.90'
$ws.Range("E23").Value = '  +2.35%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.11'
$ws.Range("E24").Value = '  +4.99%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.96'
$ws.Range("E25").Value = '  +0.80%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.76'
$ws.Range("E26").Value = '  +0.65%  '
# Row 27
$ws.Range("E27").Value = '  +0.11%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.46'
$ws.Range("E28").Value = '  -1.93%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.69'
$ws.Range("E29").Value = '  -0.43%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.15'
$ws.Range("E30").Value = '  -0.86%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.10'
$ws.Range("E31").Value = '  -0.36%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.13'
$ws.Range("E32").Value = '  +2.95%  '
# Row 33
$ws.Range("E33").Value = '  +11.67%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.19'
$ws.Range("E34").Value = '  +13.67%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0790'
$ws.Range("E35").Value = '  +5.21%  '
# Row 36
$ws.Range("E36").Value = '  +2.68%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '28.94'
$ws.Range("E37").Value = '  +9.76%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.75'
$ws.Range("E38").Value = '  +4.02%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.17'
$ws.Range("E39").Value = '  +0.34%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0321'
$ws.Range("E40").Value = '  +8.72%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.30'
$ws.Range("E41").Value = '  +5.59%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.61'
$ws.Range("E42").Value = '  +2.73%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.88'
$ws.Range("E43").Value = '  +4.50%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.15'
$ws.Range("E44").Value = '  +1.34%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.03'
$ws.Range("E45").Value = '  +0.10%  '
# Row 46
$ws.Range("E46").Value = '  +2.34%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.04'
$ws.Range("E47").Value = '  +5.28%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.103'
$ws.Range("E48").Value = '  +1.90%  '
# Row 49
$ws.Range("E49").Value = '  +1.32%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.28%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.19'
$ws.Range("E51").Value = '  +2.31%  '
